$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 7, shifting existing rows 7+ down by 2
$ws.Rows("7:8").Insert()

# The inserted rows pick up the 5-column span of their neighbours; drop the
# unused D:E cells on the two new rows so they stay genuinely blank.
$ws.Range("D7:E8").Clear()

# Inserting inside the validated C4:C8 range auto-extends the dropdown
# validation down to the shifted rows (now C4:C10); restore it to the
# original C4:C8 extent by dropping validation from the shifted rows.
$ws.Range("C9:C10").Validation.Delete()

# Populate the new row 7 (Control 2.16) - leading apostrophe forces text
# entry (matches existing "2.1"/"2.15" style Control IDs) then ClearFormats
# removes the quote-prefix formatting artifact so no new cell style is used.
$ws.Cells.Item(7, 1).Value = "'2.16"
$ws.Cells.Item(7, 1).ClearFormats()
$ws.Cells.Item(7, 2).Value = "RAG Source Integrity Validation"
$ws.Cells.Item(7, 3).Value = "Not Started"

# Populate the new row 8 (Control 2.17)
$ws.Cells.Item(8, 1).Value = "'2.17"
$ws.Cells.Item(8, 1).ClearFormats()
$ws.Cells.Item(8, 2).Value = "Multi-Agent Orchestration Limits"
$ws.Cells.Item(8, 3).Value = "Not Started"
